$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("exp_base")

# Row 25: mirai_infect_v4from - bin_download_from
$ws.Range("A25").Value = "M_N_ATK"
$ws.Range("B25").Value = "mirai_infect_v4from"
$ws.Range("C25").Value = "bin_download_from"
$ws.Range("D25").Value = "thing"
$ws.Range("E25").Value = "192.168.5.1"
$ws.Range("F25").Value = "TCP"
$ws.Range("G25").Value = "any"
$ws.Range("H25").Value = 80

$ws.Range("J25").Formula = "=CONCAT(`"{'type':'`",A25,`"',`")"
$ws.Range("K25").Formula = "=CONCAT(`"'name':'`",B25,`"-`",C25,`"',`")"
$ws.Range("L25").Formula = "=CONCAT(`"'nw_src':'`",D25,`"',`")"
$ws.Range("M25").Formula = "=CONCAT(`"'nw_dst':'`",E25,`"',`")"
$ws.Range("N25").Formula = "=CONCAT(`"'transport':'`",F25,`"',`")"
$ws.Range("O25").Formula = "=CONCAT(`"'tp_src':'`",G25,`"',`")"
$ws.Range("P25").Formula = "=CONCAT(`"'tp_dst':'`",H25,`"'}, `")"

# Row 26: mirai_infect_v4to - bin_download_to
$ws.Range("A26").Value = "M_N_ATK"
$ws.Range("B26").Value = "mirai_infect_v4to"
$ws.Range("C26").Value = "bin_download_to"
$ws.Range("D26").Value = "192.168.5.1"
$ws.Range("E26").Value = "thing"
$ws.Range("F26").Value = "TCP"
$ws.Range("G26").Value = 80
$ws.Range("H26").Value = "any"

$ws.Range("J26").Formula = "=CONCAT(`"{'type':'`",A26,`"',`")"
$ws.Range("K26").Formula = "=CONCAT(`"'name':'`",B26,`"-`",C26,`"',`")"
$ws.Range("L26").Formula = "=CONCAT(`"'nw_src':'`",D26,`"',`")"
$ws.Range("M26").Formula = "=CONCAT(`"'nw_dst':'`",E26,`"',`")"
$ws.Range("N26").Formula = "=CONCAT(`"'transport':'`",F26,`"',`")"
$ws.Range("O26").Formula = "=CONCAT(`"'tp_src':'`",G26,`"',`")"
$ws.Range("P26").Formula = "=CONCAT(`"'tp_dst':'`",H26,`"'}, `")"

$ws.Activate()
$ws.Range("P18").Select()
